$wb = $excel.ActiveWorkbook

# --- Existing sheets -------------------------------------------------
$wsSearch = $wb.Worksheets.Item(1)   # "Search"
$wsHotel  = $wb.Worksheets.Item(2)   # "Hotel"

# The Search sheet's selection moves from A1 to B2 (and it stops being the
# tab-selected sheet once Reservation becomes active below).
$wsSearch.Range("B2").Select()

# --- New "Reservation" sheet -----------------------------------------
# Copy the Hotel sheet so the new sheet inherits the same sheetFormatPr
# (defaultColWidth/defaultRowHeight/outlineLevelRow), page margins and
# namespace declarations, then strip its content and rebuild it.
$wsHotel.Copy($null, $wsHotel)
$wsRes = $wb.Worksheets.Item(3)
$wsRes.Name = "Reservation"
$wsRes.Cells.Clear()

$wsRes.Range("A1").Value = "Bed"
$wsRes.Range("B1").Value = "Amount"
$wsRes.Range("A2").Value = "Queen"
$wsRes.Range("B2").Value = 1

# Header row: centered both ways. Data row: centered horizontally.
$wsRes.Range("A1:B1").HorizontalAlignment = -4108
$wsRes.Range("A1:B1").VerticalAlignment = -4108
$wsRes.Range("A2:B2").HorizontalAlignment = -4108

# Final selection/active cell on the new sheet.
$wsRes.Range("B4").Select()
